$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster-means values (B2:V6), per commit "Ajustes para la UX/UI"
$data = @(
    @(141.3076923076923, 22.69230769230769, 42.42307692307692, 0.5413076923076923, 2.615384615384615, 8, 0.1985, 9.307692307692308, 13.73076923076923, 0.6404615384615385, 14.34615384615385, 27.84615384615385, 42.19230769230769, 10.61538461538461, 4.653846153846154, 5.307692307692307, 7.692307692307693, 16.38461538461538, 57.30769230769231, 5.461538461538462, 4.038461538461538),
    @(122.9298245614035, 17.40350877192983, 37.66666666666666, 0.4620350877192982, 7.087719298245614, 19.05263157894737, 0.3646491228070176, 7.87719298245614, 9.842105263157896, 0.8077017543859649, 4.333333333333333, 13.08771929824561, 17.42105263157895, 12.19298245614035, 4.210526315789473, 1.578947368421053, 6.333333333333333, 10.89473684210526, 49.7719298245614, 5.035087719298246, 2.035087719298246),
    @(161.5416666666667, 38.20833333333334, 73.91666666666667, 0.5199583333333333, 8, 21, 0.3650416666666667, 24.5, 32.79166666666666, 0.7587083333333333, 9.125, 36.41666666666666, 45.54166666666666, 21.79166666666667, 5.708333333333333, 5.083333333333333, 13.75, 14.41666666666667, 108.9166666666667, 4.916666666666667, 4.916666666666667),
    @(169.2448979591837, 28.79591836734694, 66.18367346938776, 0.4376530612244898, 11.63265306122449, 33, 0.3561224489795918, 12, 14.93877551020408, 0.7896734693877552, 5.163265306122449, 18.40816326530612, 23.57142857142857, 18.3265306122449, 6.122448979591836, 2.346938775510204, 8.836734693877551, 13.20408163265306, 81.22448979591837, 5.489795918367347, 4.346938775510204),
    @(186.9583333333333, 47.33333333333334, 96.79166666666667, 0.48775, 14.25, 37.95833333333334, 0.3768333333333334, 33.125, 38.5, 0.8590416666666667, 4.333333333333333, 21.79166666666667, 26.125, 30.91666666666667, 5.791666666666667, 2.458333333333333, 15.58333333333333, 13.41666666666667, 142.0416666666667, 5.5, 5.25)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
